# Auto-generated edit script: update Leve profit-tracking values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 812.0714
$ws.Range("I19").Value = 367.54544
$ws.Range("J19").Value = 1099.7059
$ws.Range("K19").Value = 367.54544
$ws.Range("L19").Value = 1099.7059
$ws.Range("M19").Value = -192.54544
$ws.Range("N19").Value = -1449.7059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2600
$ws.Range("J13").Value = 2600
$ws.Range("L13").Value = 2600
$ws.Range("N13").Value = -2888

$ws.Range("H32").Value = 31042.893
$ws.Range("I32").Value = 5441.5474
$ws.Range("J32").Value = 483333.34
$ws.Range("K32").Value = 5441.5474
$ws.Range("L32").Value = 483333.34
$ws.Range("M32").Value = -5154.5474
$ws.Range("N32").Value = -483907.34

$ws.Range("H61").Value = 1235.3125
$ws.Range("I61").Value = 869.5714
$ws.Range("J61").Value = 1519.7778
$ws.Range("K61").Value = 869.5714
$ws.Range("L61").Value = 1519.7778
$ws.Range("M61").Value = -657.5714
$ws.Range("N61").Value = -1943.7778

$ws.Range("H63").Value = 2586.0908
$ws.Range("I63").Value = 2312.25
$ws.Range("J63").Value = 3316.3333
$ws.Range("K63").Value = 2312.25
$ws.Range("L63").Value = 3316.3333
$ws.Range("M63").Value = -1626.25
$ws.Range("N63").Value = -4688.3333

$ws.Range("H66").Value = 2586.0908
$ws.Range("I66").Value = 2312.25
$ws.Range("J66").Value = 3316.3333
$ws.Range("K66").Value = 11561.25
$ws.Range("L66").Value = 16581.6665
$ws.Range("M66").Value = -8129.25
$ws.Range("N66").Value = -23445.6665

$ws.Range("H74").Value = 2781.9
$ws.Range("I74").Value = 1903.238
$ws.Range("J74").Value = 4832.1113
$ws.Range("K74").Value = 1903.238
$ws.Range("L74").Value = 4832.1113
$ws.Range("M74").Value = -1029.238
$ws.Range("N74").Value = -6580.1113

$ws.Range("H77").Value = 2781.9
$ws.Range("I77").Value = 1903.238
$ws.Range("J77").Value = 4832.1113
$ws.Range("K77").Value = 9516.190000000001
$ws.Range("L77").Value = 24160.5565
$ws.Range("M77").Value = -5148.190000000001
$ws.Range("N77").Value = -32896.5565

$ws.Range("H98").Value = 4919.4
$ws.Range("J98").Value = 4919.4
$ws.Range("L98").Value = 4919.4
$ws.Range("N98").Value = -10909.4

$ws.Range("H122").Value = 2491.389
$ws.Range("I122").Value = 3010.2856
$ws.Range("J122").Value = 2161.182
$ws.Range("K122").Value = 9030.856800000001
$ws.Range("L122").Value = 6483.545999999999
$ws.Range("M122").Value = -6580.856800000001
$ws.Range("N122").Value = -11383.546

$ws.Range("H136").Value = 1235.3125
$ws.Range("I136").Value = 869.5714
$ws.Range("J136").Value = 1519.7778
$ws.Range("K136").Value = 2608.7142
$ws.Range("L136").Value = 4559.3334
$ws.Range("M136").Value = -58.71420000000035
$ws.Range("N136").Value = -9659.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 15946.8
$ws.Range("J35").Value = 15946.8
$ws.Range("L35").Value = 15946.8
$ws.Range("N35").Value = -16566.8

$ws.Range("H82").Value = 22618.4
$ws.Range("I82").Value = 4559.3335
$ws.Range("J82").Value = 30358
$ws.Range("K82").Value = 4559.3335
$ws.Range("L82").Value = 30358
$ws.Range("M82").Value = -4176.3335
$ws.Range("N82").Value = -31124

$ws.Range("H85").Value = 22618.4
$ws.Range("I85").Value = 4559.3335
$ws.Range("J85").Value = 30358
$ws.Range("K85").Value = 4559.3335
$ws.Range("L85").Value = 30358
$ws.Range("M85").Value = -3233.3335
$ws.Range("N85").Value = -33010

$ws.Range("H134").Value = 1948.625
$ws.Range("I134").Value = 2068.9707
$ws.Range("J134").Value = 1266.6666
$ws.Range("K134").Value = 6206.9121
$ws.Range("L134").Value = 3799.9998
$ws.Range("M134").Value = -3671.9121
$ws.Range("N134").Value = -8869.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19235.316
$ws.Range("I31").Value = 23346.979
$ws.Range("J31").Value = 3816.5833
$ws.Range("K31").Value = 23346.979
$ws.Range("L31").Value = 3816.5833
$ws.Range("M31").Value = -23051.979
$ws.Range("N31").Value = -4406.5833

$ws.Range("H34").Value = 19235.316
$ws.Range("I34").Value = 23346.979
$ws.Range("J34").Value = 3816.5833
$ws.Range("K34").Value = 23346.979
$ws.Range("L34").Value = 3816.5833
$ws.Range("M34").Value = -23144.979
$ws.Range("N34").Value = -4220.5833

$ws.Range("H94").Value = 1247.3334
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1247.3334
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1247.3334
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2149.3334

$ws.Range("H99").Value = 22798.6
$ws.Range("J99").Value = 27253.25
$ws.Range("L99").Value = 27253.25
$ws.Range("N99").Value = -30249.25

$ws.Range("H124").Value = 41990
$ws.Range("J124").Value = 41990
$ws.Range("L124").Value = 41990
$ws.Range("N124").Value = -46900

$ws.Range("H126").Value = 22798.6
$ws.Range("J126").Value = 27253.25
$ws.Range("L126").Value = 81759.75
$ws.Range("N126").Value = -86699.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3092.0588
$ws.Range("I132").Value = 2468.9285
$ws.Range("K132").Value = 22220.3565
$ws.Range("M132").Value = -19690.3565

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 12500700
$ws.Range("I3").Value = 12500700
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 12500700
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -12500584
$ws.Range("N3").ClearContents()

$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338

$ws.Range("H12").Value = 4732435
$ws.Range("I12").Value = 6834111
$ws.Range("J12").Value = 3663
$ws.Range("K12").Value = 6834111
$ws.Range("L12").Value = 3663
$ws.Range("M12").Value = -6833971
$ws.Range("N12").Value = -3943

$ws.Range("H132").Value = 2354.2188
$ws.Range("I132").Value = 1767.762
$ws.Range("J132").Value = 3473.818
$ws.Range("K132").Value = 5303.286
$ws.Range("L132").Value = 10421.454
$ws.Range("M132").Value = -2773.286
$ws.Range("N132").Value = -15481.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 53334.668
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 53334.668
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 53334.668
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -53678.668

$ws.Range("H40").Value = 85437.336
$ws.Range("I40").Value = 501000
$ws.Range("K40").Value = 501000
$ws.Range("M40").Value = -500864

$ws.Range("H81").Value = 22827.75
$ws.Range("J81").Value = 22827.75
$ws.Range("L81").Value = 22827.75
$ws.Range("N81").Value = -24823.75

$ws.Range("H84").Value = 22827.75
$ws.Range("J84").Value = 22827.75
$ws.Range("L84").Value = 68483.25
$ws.Range("N84").Value = -78467.25

$ws.Range("H119").Value = 43990
$ws.Range("J119").Value = 43990
$ws.Range("L119").Value = 43990
$ws.Range("N119").Value = -53666

$ws.Range("H122").Value = 2602.7222
$ws.Range("I122").Value = 1943.2858
$ws.Range("J122").Value = 3022.3635
$ws.Range("K122").Value = 5829.857400000001
$ws.Range("L122").Value = 9067.0905
$ws.Range("M122").Value = -3379.857400000001
$ws.Range("N122").Value = -13967.0905

$ws.Range("H132").Value = 3284.5186
$ws.Range("I132").Value = 3159.9565
$ws.Range("K132").Value = 9479.869499999999
$ws.Range("M132").Value = -6949.869499999999

$ws.Range("H136").Value = 1942
$ws.Range("I136").Value = 1999.8334
$ws.Range("J136").Value = 1884.1666
$ws.Range("K136").Value = 5999.5002
$ws.Range("L136").Value = 5652.4998
$ws.Range("M136").Value = -3449.5002
$ws.Range("N136").Value = -10752.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 28445
$ws.Range("J124").Value = 28445
$ws.Range("L124").Value = 28445
$ws.Range("N124").Value = -38265

$ws.Range("H136").Value = 1004.1177
$ws.Range("I136").Value = 760.9091
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 2282.7273
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = 267.2727
$ws.Range("N136").Value = -9450

$ws.Range("H140").Value = 53413.355
$ws.Range("J140").Value = 53413.355
$ws.Range("L140").Value = 53413.355
$ws.Range("N140").Value = -63773.355
